{"js": "// Update the date paragraph and the 25 arithmetic-answer cells in the\n// single practice table. Each table row/cell is addressed positionally\n// (row, column) so duplicate old values (e.g. \"433\u00d72=866\" appearing\n// twice) are replaced independently with their own new values.\n\n// 1) Date heading paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.getRange().insertText(\"2025-10-10 Friday\", Word.InsertLocation.replace);\n\n// 2) The answers table. Only 5 of the 20 rows contain text (rows 0, 4, 9,\n// 14, 19); each of those rows has 5 cells, for 25 total replacements.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> [newCol0, newCol1, newCol2, newCol3, newCol4]\nconst rowUpdates = {\n  0: [\"437\u00d72=874\", \"639\u00d73=1917\", \"729\u00d78=5832\", \"626\u00d78=5008\", \"222\u00d72=444\"],\n  4: [\"548\u00d72=1096\", \"597\u00d76=3582\", \"788\u00d75=3940\", \"342\u00d79=3078\", \"772\u00d73=2316\"],\n  9: [\"734\u00d76=4404\", \"225\u00d73=675\", \"741\u00d74=2964\", \"621\u00d74=2484\", \"383\u00d75=1915\"],\n  14: [\"967\u00d76=5802\", \"306\u00d79=2754\", \"554\u00d72=1108\", \"261\u00d76=1566\", \"984\u00d77=6888\"],\n  19: [\"904\u00d77=6328\", \"445\u00d76=2670\", \"511\u00d73=1533\", \"594\u00d79=5346\", \"473\u00d72=946\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newValues = rowUpdates[rowIndex];\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const para = cell.body.paragraphs.getFirst();\n    para.getRange().insertText(newValues[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 arithmetic-answer cells in the\n# single practice table. Each table cell is addressed positionally\n# (row, column) so duplicate old values (e.g. \"433x2=866\" appearing\n# twice) are replaced independently with their own new values.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph (first paragraph of the body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-10-10 Friday\"\n\n# 2) The answers table. Only 5 of the 20 rows contain text (rows 1, 5, 10,\n# 15, 20 in 1-based COM indexing); each of those rows has 5 cells, for 25\n# total replacements.\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"437\u00d72=874\", \"639\u00d73=1917\", \"729\u00d78=5832\", \"626\u00d78=5008\", \"222\u00d72=444\")\n    5  = @(\"548\u00d72=1096\", \"597\u00d76=3582\", \"788\u00d75=3940\", \"342\u00d79=3078\", \"772\u00d73=2316\")\n    10 = @(\"734\u00d76=4404\", \"225\u00d73=675\", \"741\u00d74=2964\", \"621\u00d74=2484\", \"383\u00d75=1915\")\n    15 = @(\"967\u00d76=5802\", \"306\u00d79=2754\", \"554\u00d72=1108\", \"261\u00d76=1566\", \"984\u00d77=6888\")\n    20 = @(\"904\u00d77=6328\", \"445\u00d76=2670\", \"511\u00d73=1533\", \"594\u00d79=5346\", \"473\u00d72=946\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $newValues[$col - 1]\n    }\n}\n"}
